# Automatic edit reproduced in Card5 (and the companion clear on Card6),
# matching the commit "تعديل تلقائي في شيت Card5 by admin at 2025-12-06 18:33:23".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Card6: every data cell in D2:L12 currently holds the placeholder text
# "nan" (no real data was ever entered). The edit clears all of them back
# to blank cells.
# ---------------------------------------------------------------------
$wsCard6 = $wb.Worksheets.Item("Card6")
for ($r = 2; $r -le 12; $r++) {
    for ($c = 4; $c -le 12; $c++) {
        $wsCard6.Cells.Item($r, $c).ClearContents()
    }
}

# ---------------------------------------------------------------------
# Card5: the "card" id in A2 was mis-entered as 2 and should be 5 (it
# matches the card-5 ranges used elsewhere on this sheet, e.g. row 8).
# Use a quote-prefixed formula so Excel keeps storing it as text (the
# rest of the sheet stores every value, numeric-looking or not, as text).
# ---------------------------------------------------------------------
$wsCard5 = $wb.Worksheets.Item("Card5")
$wsCard5.Cells.Item(2, 1).Formula = "'5"

# Every still-blank cell in D2:O13 gets back-filled with the placeholder
# text "nan" (mirrors the "nan" placeholders already used on this sheet,
# e.g. D3, D6, D7 ...), while any cell that already holds real data is
# left untouched.
for ($r = 2; $r -le 13; $r++) {
    for ($c = 4; $c -le 15; $c++) {
        $cell = $wsCard5.Cells.Item($r, $c)
        if ($cell.Text -eq "") {
            $cell.Value = "nan"
        }
    }
}
